# Fruta / hortaliza, semanal
# Insert 3 new price records (August Red, Región de O'Higgins, date 2022-03-17)
# just above the existing "Nectar Crest" row (row 291), pushing the rest of the
# table down by 3 rows (old A1:T343 -> new A1:T346).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 291..343 down by inserting 3 blank rows at 291.
$ws.Range("A291:A293").EntireRow.Insert()

# Common (constant across this sub-block) column values.
$marketId   = 8
$market     = "Terminal La Palmera de La Serena"
$region     = "Coquimbo"
$codreg     = 4
$tipo       = "Fruta"
$productoId = 100103
$producto   = "Frutos de hueso (carozo)"
$categoriaId = 100103006
$categoria  = "Nectarín"
$unidad     = "$/bins (420 kilos)"
$origen     = "Región de O'Higgins"

# New row 291: August Red, Especial
$ws.Cells.Item(291, 1).Value  = $marketId
$ws.Cells.Item(291, 2).Value  = $market
$ws.Cells.Item(291, 3).Value  = $region
$ws.Cells.Item(291, 4).Value  = 44637
$ws.Cells.Item(291, 5).Value  = $codreg
$ws.Cells.Item(291, 6).Value  = $tipo
$ws.Cells.Item(291, 7).Value  = $productoId
$ws.Cells.Item(291, 8).Value  = $producto
$ws.Cells.Item(291, 9).Value  = $categoriaId
$ws.Cells.Item(291, 10).Value = $categoria
$ws.Cells.Item(291, 11).Value = "August Red"
$ws.Cells.Item(291, 12).Value = "Especial"
$ws.Cells.Item(291, 13).Value = 16
$ws.Cells.Item(291, 14).Value = 450000
$ws.Cells.Item(291, 15).Value = 460000
$ws.Cells.Item(291, 16).Value = 455000
$ws.Cells.Item(291, 17).Value = $unidad
$ws.Cells.Item(291, 18).Value = $origen
$ws.Cells.Item(291, 19).Value = 1083
$ws.Cells.Item(291, 20).Value = 420

# New row 292: August Red, Primera
$ws.Cells.Item(292, 1).Value  = $marketId
$ws.Cells.Item(292, 2).Value  = $market
$ws.Cells.Item(292, 3).Value  = $region
$ws.Cells.Item(292, 4).Value  = 44637
$ws.Cells.Item(292, 5).Value  = $codreg
$ws.Cells.Item(292, 6).Value  = $tipo
$ws.Cells.Item(292, 7).Value  = $productoId
$ws.Cells.Item(292, 8).Value  = $producto
$ws.Cells.Item(292, 9).Value  = $categoriaId
$ws.Cells.Item(292, 10).Value = $categoria
$ws.Cells.Item(292, 11).Value = "August Red"
$ws.Cells.Item(292, 12).Value = "Primera"
$ws.Cells.Item(292, 13).Value = 20
$ws.Cells.Item(292, 14).Value = 420000
$ws.Cells.Item(292, 15).Value = 430000
$ws.Cells.Item(292, 16).Value = 425000
$ws.Cells.Item(292, 17).Value = $unidad
$ws.Cells.Item(292, 18).Value = $origen
$ws.Cells.Item(292, 19).Value = 1012
$ws.Cells.Item(292, 20).Value = 420

# New row 293: August Red, Segunda
$ws.Cells.Item(293, 1).Value  = $marketId
$ws.Cells.Item(293, 2).Value  = $market
$ws.Cells.Item(293, 3).Value  = $region
$ws.Cells.Item(293, 4).Value  = 44637
$ws.Cells.Item(293, 5).Value  = $codreg
$ws.Cells.Item(293, 6).Value  = $tipo
$ws.Cells.Item(293, 7).Value  = $productoId
$ws.Cells.Item(293, 8).Value  = $producto
$ws.Cells.Item(293, 9).Value  = $categoriaId
$ws.Cells.Item(293, 10).Value = $categoria
$ws.Cells.Item(293, 11).Value = "August Red"
$ws.Cells.Item(293, 12).Value = "Segunda"
$ws.Cells.Item(293, 13).Value = 20
$ws.Cells.Item(293, 14).Value = 380000
$ws.Cells.Item(293, 15).Value = 390000
$ws.Cells.Item(293, 16).Value = 385000
$ws.Cells.Item(293, 17).Value = $unidad
$ws.Cells.Item(293, 18).Value = $origen
$ws.Cells.Item(293, 19).Value = 917
$ws.Cells.Item(293, 20).Value = 420
